# Reverse the order of the comma-separated "Recorded By" names in column G
# for every data row of the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(",")
        if ($parts.Length -gt 1) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }
            $reversed = $trimmed[-1..-($trimmed.Length)]
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value2 = $newVal
        }
    }
}
